# Auto-generated Excel COM-interop edit script
# Applies refreshed market-price / profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) to each class leve-profit sheet, matching the
# scheduled market-data runner output.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$changes = @(
    @(2, "H", 17666.666),
    @(2, "J", 25150),
    @(2, "L", 25150),
    @(2, "N", -25376),
    @(15, "H", 1097.7646),
    @(15, "I", 1097.7646),
    @(15, "K", 3293.2938),
    @(15, "M", -3124.2938),
    @(74, "H", 0),
    @(74, "I", 0),
    @(74, "K", 0),
    @(74, "M", $null),
    @(77, "H", 0),
    @(77, "I", 0),
    @(77, "K", 0),
    @(77, "M", $null),
    @(92, "H", 1643.8572),
    @(92, "I", 1663.5385),
    @(92, "K", 1663.5385),
    @(92, "M", -415.5385000000001),
    @(103, "H", 1897),
    @(103, "I", 1897),
    @(103, "K", 5691),
    @(103, "M", -5105),
    @(125, "H", 7831.3335),
    @(125, "I", 7748.5),
    @(125, "K", 69736.5),
    @(125, "M", -67276.5),
    @(135, "H", 3766.1),
    @(135, "I", 4775.143),
    @(135, "K", 42976.287),
    @(135, "M", -40441.287),
    @(138, "H", 6978.8857),
    @(138, "J", 7609.6665),
    @(138, "L", 22828.9995),
    @(138, "N", -33108.99950000001)
)
foreach ($chg in $changes) {
    $cellRef = [string]$chg[1] + [string]$chg[0]
    $ws.Range($cellRef).Value = $chg[2]
}

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$changes = @(
    @(2, "H", 999),
    @(2, "I", 998.5),
    @(2, "K", 998.5),
    @(2, "M", -885.5),
    @(32, "H", 3114.3044),
    @(32, "I", 2514.5476),
    @(32, "K", 2514.5476),
    @(32, "M", -2227.5476),
    @(45, "H", 3003.9167),
    @(45, "I", 2689.5),
    @(45, "K", 2689.5),
    @(45, "M", -2312.5),
    @(110, "H", 3288.9333),
    @(110, "I", 3310.2144),
    @(110, "J", 2991),
    @(110, "K", 3310.2144),
    @(110, "L", 2991),
    @(110, "M", -1265.2144),
    @(110, "N", -7081),
    @(116, "H", 999),
    @(116, "I", 998.5),
    @(116, "K", 998.5),
    @(116, "M", 1295.5),
    @(122, "H", 1570.3684),
    @(122, "I", 1601.2778),
    @(122, "J", 1014),
    @(122, "K", 4803.8334),
    @(122, "L", 3042),
    @(122, "M", -2353.8334),
    @(122, "N", -7942),
    @(132, "H", 3225.087),
    @(132, "I", 2629.5386),
    @(132, "J", 3999.3),
    @(132, "K", 7888.6158),
    @(132, "L", 11997.9),
    @(132, "M", -5358.6158),
    @(132, "N", -17057.9)
)
foreach ($chg in $changes) {
    $cellRef = [string]$chg[1] + [string]$chg[0]
    $ws.Range($cellRef).Value = $chg[2]
}

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$changes = @(
    @(3, "H", 999),
    @(3, "I", 998.5),
    @(3, "K", 998.5),
    @(3, "M", -884.5),
    @(20, "H", 6332.6665),
    @(20, "I", 6332.6665),
    @(20, "K", 6332.6665),
    @(20, "M", -6085.6665),
    @(86, "H", 1713.2858),
    @(86, "I", 1498.6),
    @(86, "K", 1498.6),
    @(86, "M", -375.5999999999999),
    @(89, "H", 1713.2858),
    @(89, "I", 1498.6),
    @(89, "K", 7493),
    @(89, "M", -1877),
    @(94, "H", 1051.5),
    @(94, "I", 903),
    @(94, "K", 903),
    @(94, "M", -452),
    @(105, "H", 2854.5715),
    @(105, "I", 2854.5715),
    @(105, "K", 2854.5715),
    @(105, "M", -1107.5715),
    @(133, "H", 75000),
    @(133, "J", 75000),
    @(133, "L", 75000),
    @(133, "N", -85120)
)
foreach ($chg in $changes) {
    $cellRef = [string]$chg[1] + [string]$chg[0]
    $ws.Range($cellRef).Value = $chg[2]
}

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$changes = @(
    @(16, "H", 8389),
    @(16, "I", 8389),
    @(16, "K", 8389),
    @(16, "M", -8102),
    @(31, "H", 5529.905),
    @(31, "I", 3271.5),
    @(31, "J", 8541.111000000001),
    @(31, "K", 3271.5),
    @(31, "L", 8541.111000000001),
    @(31, "M", -2976.5),
    @(31, "N", -9131.111000000001),
    @(34, "H", 5529.905),
    @(34, "I", 3271.5),
    @(34, "J", 8541.111000000001),
    @(34, "K", 3271.5),
    @(34, "L", 8541.111000000001),
    @(34, "M", -3069.5),
    @(34, "N", -8945.111000000001),
    @(58, "H", 2210.0908),
    @(58, "I", 2192.375),
    @(58, "J", 2257.3333),
    @(58, "K", 2192.375),
    @(58, "L", 2257.3333),
    @(58, "M", -1989.375),
    @(58, "N", -2663.3333),
    @(107, "H", 1323.1428),
    @(107, "I", 902.1667),
    @(107, "K", 902.1667),
    @(107, "M", 1017.8333),
    @(113, "H", 8389),
    @(113, "I", 8389),
    @(113, "K", 8389),
    @(113, "M", -6219),
    @(136, "H", 2210.0908),
    @(136, "I", 2192.375),
    @(136, "J", 2257.3333),
    @(136, "K", 6577.125),
    @(136, "L", 6771.999899999999),
    @(136, "M", -4027.125),
    @(136, "N", -11871.9999)
)
foreach ($chg in $changes) {
    $cellRef = [string]$chg[1] + [string]$chg[0]
    $ws.Range($cellRef).Value = $chg[2]
}

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$changes = @(
    @(12, "H", 295.9091),
    @(12, "I", 219.8),
    @(12, "J", 359.33334),
    @(12, "K", 659.4000000000001),
    @(12, "L", 1078.00002),
    @(12, "M", -486.4000000000001),
    @(12, "N", -1424.00002)
)
foreach ($chg in $changes) {
    $cellRef = [string]$chg[1] + [string]$chg[0]
    $ws.Range($cellRef).Value = $chg[2]
}

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$changes = @(
    @(14, "H", 10000),
    @(14, "I", 0),
    @(14, "J", 10000),
    @(14, "K", 0),
    @(14, "L", 10000),
    @(14, "M", $null),
    @(14, "N", -10336),
    @(36, "H", 2017),
    @(36, "I", 2017),
    @(36, "K", 2017),
    @(36, "M", -1532),
    @(43, "H", 12414.571),
    @(43, "I", 2017),
    @(43, "J", 14147.5),
    @(43, "K", 2017),
    @(43, "L", 14147.5),
    @(43, "M", -1866),
    @(43, "N", -14449.5),
    @(70, "H", 3499),
    @(70, "I", 3499),
    @(70, "K", 3499),
    @(70, "M", -3229),
    @(73, "H", 3499),
    @(73, "I", 3499),
    @(73, "K", 3499),
    @(73, "M", -2563),
    @(80, "H", 3860.7),
    @(80, "I", 3023.5833),
    @(80, "J", 5116.375),
    @(80, "K", 3023.5833),
    @(80, "L", 5116.375),
    @(80, "M", -2025.5833),
    @(80, "N", -7112.375),
    @(83, "H", 3860.7),
    @(83, "I", 3023.5833),
    @(83, "J", 5116.375),
    @(83, "K", 15117.9165),
    @(83, "L", 25581.875),
    @(83, "M", -10125.9165),
    @(83, "N", -35565.875),
    @(102, "H", 1542.8),
    @(102, "I", 1612.6666),
    @(102, "J", 914),
    @(102, "K", 1612.6666),
    @(102, "L", 914),
    @(102, "M", 9.333399999999983),
    @(102, "N", -4158),
    @(122, "H", 622.5),
    @(122, "I", 496.66666),
    @(122, "J", 1000),
    @(122, "K", 1489.99998),
    @(122, "L", 3000),
    @(122, "M", 960.0000199999999),
    @(122, "N", -7900),
    @(126, "H", 2030.4286),
    @(126, "I", 1737.6666),
    @(126, "K", 5212.9998),
    @(126, "M", -2742.9998),
    @(132, "H", 4977.0454),
    @(132, "I", 4913.0835),
    @(132, "J", 5053.8),
    @(132, "K", 14739.2505),
    @(132, "L", 15161.4),
    @(132, "M", -12209.2505),
    @(132, "N", -20221.4),
    @(134, "H", 0),
    @(134, "J", 0),
    @(134, "L", 0),
    @(134, "N", $null)
)
foreach ($chg in $changes) {
    $cellRef = [string]$chg[1] + [string]$chg[0]
    $ws.Range($cellRef).Value = $chg[2]
}

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$changes = @(
    @(7, "H", 3399.4),
    @(7, "I", 2999.3333),
    @(7, "J", 3999.5),
    @(7, "K", 2999.3333),
    @(7, "L", 3999.5),
    @(7, "M", -2887.3333),
    @(7, "N", -4223.5),
    @(40, "H", 3021.7778),
    @(40, "I", 3042.4285),
    @(40, "J", 2949.5),
    @(40, "K", 3042.4285),
    @(40, "L", 2949.5),
    @(40, "M", -2906.4285),
    @(40, "N", -3221.5),
    @(122, "H", 3960),
    @(122, "I", 3960),
    @(122, "K", 11880),
    @(122, "M", -9430),
    @(126, "H", 3399.4),
    @(126, "I", 2999.3333),
    @(126, "J", 3999.5),
    @(126, "K", 8997.999899999999),
    @(126, "L", 11998.5),
    @(126, "M", -6527.999899999999),
    @(126, "N", -16938.5),
    @(136, "H", 25526.428),
    @(136, "I", 3008.125),
    @(136, "K", 9024.375),
    @(136, "M", -6474.375)
)
foreach ($chg in $changes) {
    $cellRef = [string]$chg[1] + [string]$chg[0]
    $ws.Range($cellRef).Value = $chg[2]
}

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$changes = @(
    @(107, "H", 342),
    @(107, "I", 310.8),
    @(107, "K", 932.4000000000001),
    @(107, "M", 987.5999999999999),
    @(126, "H", 2017.5),
    @(126, "I", 1856.8334),
    @(126, "J", 2499.5),
    @(126, "K", 5570.5002),
    @(126, "L", 7498.5),
    @(126, "M", -3100.5002),
    @(126, "N", -12438.5),
    @(132, "H", 2715.7812),
    @(132, "I", 2424.92),
    @(132, "K", 7274.76),
    @(132, "M", -4744.76)
)
foreach ($chg in $changes) {
    $cellRef = [string]$chg[1] + [string]$chg[0]
    $ws.Range($cellRef).Value = $chg[2]
}
